$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column I: QUANTITY ALERT
$ws.Range("I1").Value = "QUANTITY ALERT"

# Populate QUANTITY ALERT = 500 for all data rows except row 4
$rows = 2..22
foreach ($r in $rows) {
    if ($r -eq 4) { continue }
    $ws.Cells.Item($r, 9).Value = 500
}

# Match the final active selection from the diff
$ws.Range("I5").Select()
